$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and B28/B29, C28/C29 swap)

# Row 2
$ws.Range("D2").Value = "'27.688.69"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3
$ws.Range("D3").Value = "'1.639.40"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
$ws.Range("D4").Value = "'1.00"

# Row 5
$ws.Range("D5").Value = "'212.36"
$ws.Range("E5").Value = "  -0.13%  "

# Row 6
$ws.Range("D6").Value = "'0.522"
$ws.Range("E6").Value = "  -1.51%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'23.08"
$ws.Range("E8").Value = "  -2.25%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("E10").Value = "  -0.13%  "

# Row 11
$ws.Range("D11").Value = "'0.0893"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12
$ws.Range("D12").Value = "'1.871.81"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13
$ws.Range("D13").Value = "'1.645.04"
$ws.Range("E13").Value = "  -0.29%  "

# Row 14
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("E15").Value = "  -5.77%  "

# Row 16
$ws.Range("D16").Value = "'64.65"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("D17").Value = "'27.669.45"
$ws.Range("E17").Value = "  +0.55%  "

# Row 18
$ws.Range("D18").Value = "'230.51"
$ws.Range("E18").Value = "  -0.66%  "

# Row 19
$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +2.18%  "

# Row 20
$ws.Range("E20").Value = "  -0.40%  "

# Row 21
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("E22").Value = "  -0.76%  "

# Row 23
$ws.Range("D23").Value = "'10.23"
$ws.Range("E23").Value = "  +4.37%  "

# Row 24
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +0.92%  "

# Row 25
$ws.Range("D25").Value = "'151.27"
$ws.Range("E25").Value = "  +1.76%  "

# Row 26
$ws.Range("E26").Value = "  -1.18%  "

# Row 27
$ws.Range("D27").Value = "'0.112"
$ws.Range("E27").Value = "  -1.60%  "

# Row 28
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'15.61"

# Row 30
$ws.Range("E30").Value = "  +0.01%  "

# Row 31
$ws.Range("D31").Value = "'0.0487"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("E32").Value = "  -0.42%  "

# Row 33
$ws.Range("D33").Value = "'1.458.53"
$ws.Range("E33").Value = "  +2.37%  "

# Row 34
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  -1.28%  "

# Row 35
$ws.Range("E35").Value = "  -1.37%  "

# Row 36
$ws.Range("E36").Value = "  -0.32%  "

# Row 37
$ws.Range("E37").Value = "  -0.73%  "

# Row 38
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  -1.23%  "

# Row 39
$ws.Range("E39").Value = "  +0.24%  "

# Row 40
$ws.Range("E40").Value = "  +9.50%  "

# Row 41
$ws.Range("D41").Value = "'69.88"
$ws.Range("E41").Value = "  +7.25%  "

# Row 42
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  -0.98%  "

# Row 44
$ws.Range("E44").Value = "  +1.43%  "

# Row 45
$ws.Range("E45").Value = "  -0.58%  "

# Row 46
$ws.Range("D46").Value = "'2.24"
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$ws.Range("D47").Value = "'1.781.49"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48
$ws.Range("E48").Value = "  +3.38%  "

# Row 49
$ws.Range("D49").Value = "'86.86"
$ws.Range("E49").Value = "  -1.47%  "

# Row 50
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("E51").Value = "  -0.15%  "

